$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Date values in column C shift forward by 21 days (re-run of the assay).
#    C2:C6  -> 45455 (was 45434)
#    C7:C46 -> 45458 (was 45437)
# ---------------------------------------------------------------------------
for ($r = 2; $r -le 6; $r++) {
  $ws.Cells.Item($r, 3).Value = 45455
}
for ($r = 7; $r -le 46; $r++) {
  $ws.Cells.Item($r, 3).Value = 45458
}

# ---------------------------------------------------------------------------
# 2. Re-style column C: drop the bordered/shaded header-style formatting and
#    use a plain date number format instead. This clears font/fill/border on
#    the whole column (incl. header cell C1) and reapplies a clean date
#    format to the data cells C2:C46.
# ---------------------------------------------------------------------------
$ws.Columns.Item(3).ClearFormats()

$ws.Range("C2").NumberFormat = "mm-dd-yy"
$ws.Range("C2").Copy()
$ws.Range("C3:C46").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Column C narrows slightly now that it carries the default (unstyled) font.
$ws.Columns.Item(3).ColumnWidth = 8.5

# ---------------------------------------------------------------------------
# 3. Cosmetic: update the saved cell selection on the sheet.
# ---------------------------------------------------------------------------
[void]$ws.Range("J12").Select()
